$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.10561669283312
$ws.Range("G2").Value = 0.00506455245660282
$ws.Range("H2").Value = 0.0926064804802312
$ws.Range("I2").Value = 0.0921262988525608
$ws.Range("L2").Value = 0.0555111539218084
$ws.Range("M2").Value = 0.13515351676554
$ws.Range("N2").Value = 0.0738321659337642
$ws.Range("O2").Value = 0.109889836172809
$ws.Range("Q2").Value = 0.00321310532969494
